$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 32   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# --- Data table updates ---
# Constants used to paste formats without minting new style records:
#   xlPasteFormats = -4122 ; xlPasteValues = -4163
$xlPasteFormats = -4122
$xlPasteValues = -4163

$ws.Range("C23").Copy() | Out-Null
$ws.Range("D14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D14").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("E23").Copy() | Out-Null
$ws.Range("E14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E14").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D15").Value = 1
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E15").Value = -100
$ws.Range("G14").Copy() | Out-Null
$ws.Range("G15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G15").Value = 1
$ws.Range("H14").Copy() | Out-Null
$ws.Range("H15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = 30
$ws.Range("M15").Value = 85.714285714285
$ws.Range("C23").Copy() | Out-Null
$ws.Range("C16").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C16").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D16").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D16").Value = 2
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E16").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E16").Value = -100
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -33.333333333333
$ws.Range("J16").Value = 70
$ws.Range("K16").Value = -18.571428571428
$ws.Range("L16").Value = -47.222222222222
$ws.Range("N16").Value = -90.086956521739
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 12
$ws.Range("I17").Value = 82
$ws.Range("J17").Value = 109
$ws.Range("K17").Value = -24.770642201834
$ws.Range("L17").Value = -16.326530612244
$ws.Range("M17").Value = 64
$ws.Range("N17").Value = -43.835616438356
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 94
$ws.Range("J18").Value = 87
$ws.Range("K18").Value = 8.045977011494
$ws.Range("L18").Value = -19.658119658119
$ws.Range("M18").Value = 8.045977011494
$ws.Range("N18").Value = -91.637010676156
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 12.5
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = 34.210526315789
$ws.Range("I19").Value = 523
$ws.Range("J19").Value = 537
$ws.Range("K19").Value = -2.607076350093
$ws.Range("L19").Value = -16.984126984127
$ws.Range("M19").Value = -22.518518518518
$ws.Range("N19").Value = -75
$ws.Range("C20").Value = 1
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D20").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D20").Value = 1
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E20").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 27
$ws.Range("J20").Value = 24
$ws.Range("K20").Value = 12.5
$ws.Range("L20").Value = -47.058823529411
$ws.Range("M20").Value = -10
$ws.Range("N20").Value = -95.368782161235
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = -13.333333333333
$ws.Range("F21").Value = 76
$ws.Range("G21").Value = 68
$ws.Range("H21").Value = 11.764705882352
$ws.Range("I21").Value = 796
$ws.Range("J21").Value = 841
$ws.Range("K21").Value = -5.350772889417
$ws.Range("L21").Value = -21.266073194856
$ws.Range("M21").Value = -12.044198895027
$ws.Range("N21").Value = -82.439885285682
$ws.Range("C22").Value = 1
$ws.Range("C23").Copy() | Out-Null
$ws.Range("D22").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D22").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("E23").Copy() | Out-Null
$ws.Range("E22").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E22").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 25
$ws.Range("K22").Value = -28.571428571428
$ws.Range("L22").Value = -30.555555555555
$ws.Range("M22").Value = -21.875
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 66
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = 10
$ws.Range("I24").Value = 891
$ws.Range("J24").Value = 854
$ws.Range("K24").Value = 4.332552693208
$ws.Range("L24").Value = -9.081632653061
$ws.Range("M24").Value = 55.226480836236
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = 7.142857142857
$ws.Range("I25").Value = 664
$ws.Range("J25").Value = 660
$ws.Range("K25").Value = 0.60606060606
$ws.Range("L25").Value = -17.412935323383
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 11
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = -35.294117647058
$ws.Range("I26").Value = 175
$ws.Range("J26").Value = 209
$ws.Range("K26").Value = -16.267942583732
$ws.Range("L26").Value = -23.5807860262
$ws.Range("M26").Value = -18.22429906542
$ws.Range("G14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D27").Value = 1
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E27").Value = -100
$ws.Range("G14").Copy() | Out-Null
$ws.Range("G27").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G27").Value = 1
$ws.Range("H14").Copy() | Out-Null
$ws.Range("H27").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 12
$ws.Range("K27").Value = 8.333333333333
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 51
$ws.Range("K28").Value = -9.803921568627
$ws.Range("L31").Value = -18.181818181818

$excel.CutCopyMode = $false
Write-Output "done"
